# Tukey HSD rejected table: add 9v to Vin, and recalibrate all stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reject" column (old column H). This shifts the old "H0"
# column (I) left into H, matching the new header row group1..p-value,H0.
$ws.Columns.Item(8).EntireColumn.Delete()

# ---- Row 2 ----
$ws.Cells.Item(2,1).Value2 = "0.5v"
$ws.Cells.Item(2,2).Value2 = "18v"
$ws.Cells.Item(2,3).Value2 = 0.01754990930677841
$ws.Cells.Item(2,4).Value2 = 0.001247209584349981
$ws.Cells.Item(2,5).Value2 = 0.03385260902920684
$ws.Cells.Item(2,6).Value2 = 4.723470142111323
$ws.Cells.Item(2,7).Value2 = 0.02379130724204015
$ws.Cells.Item(2,8).Value2 = "'False"

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value2 = "3v"
$ws.Cells.Item(3,2).Value2 = "18v"
$ws.Cells.Item(3,3).Value2 = 0.01754990930677841
$ws.Cells.Item(3,4).Value2 = 0.00123683516945973
$ws.Cells.Item(3,5).Value2 = 0.03386298344409709
$ws.Cells.Item(3,6).Value2 = 4.720466217862374
$ws.Cells.Item(3,7).Value2 = 0.02395482209035593
$ws.Cells.Item(3,8).Value2 = "'False"

# ---- Row 4 ----
$ws.Cells.Item(4,1).Value2 = "6v"
$ws.Cells.Item(4,2).Value2 = "24.5v"
$ws.Cells.Item(4,3).Value2 = 0.01870862978179839
$ws.Cells.Item(4,4).Value2 = 0.004804548390456775
$ws.Cells.Item(4,5).Value2 = 0.03261271117314001
$ws.Cells.Item(4,6).Value2 = 5.903988555535277
$ws.Cells.Item(4,7).Value2 = 0.001005322592465663
$ws.Cells.Item(4,8).Value2 = "'False"

# ---- Row 5 ----
$ws.Cells.Item(5,1).Value2 = "9v"
$ws.Cells.Item(5,2).Value2 = "18v"
$ws.Cells.Item(5,3).Value2 = 0.01407768708455626
$ws.Cells.Item(5,4).Value2 = 0.0002240114764534996
$ws.Cells.Item(5,5).Value2 = 0.02793136269265901
$ws.Cells.Item(5,6).Value2 = 4.458739717283347
$ws.Cells.Item(5,7).Value2 = 0.04305586947443718
$ws.Cells.Item(5,8).Value2 = "'False"

# ---- Row 6 ----
$ws.Cells.Item(6,1).Value2 = "12v"
$ws.Cells.Item(6,2).Value2 = "18v"
$ws.Cells.Item(6,3).Value2 = 0.01925445476132383
$ws.Cells.Item(6,4).Value2 = 0.005334826344308272
$ws.Cells.Item(6,5).Value2 = 0.0331740831783394
$ws.Cells.Item(6,6).Value2 = 6.069451010385574
$ws.Cells.Item(6,7).Value2 = 0.001
$ws.Cells.Item(6,8).Value2 = "'False"

# ---- Row 7 ----
$ws.Cells.Item(7,1).Value2 = "15v"
$ws.Cells.Item(7,2).Value2 = "24.5v"
$ws.Cells.Item(7,3).Value2 = 0.01930807137229723
$ws.Cells.Item(7,4).Value2 = 0.005660638378136242
$ws.Cells.Item(7,5).Value2 = 0.03295550436645823
$ws.Cells.Item(7,6).Value2 = 6.207743333265988
$ws.Cells.Item(7,7).Value2 = 0.001
$ws.Cells.Item(7,8).Value2 = "'False"

# ---- Row 8 ----
$ws.Cells.Item(8,1).Value2 = "18v"
$ws.Cells.Item(8,2).Value2 = "24.5v"
$ws.Cells.Item(8,3).Value2 = 0.02696767826187999
$ws.Cells.Item(8,4).Value2 = 0.01348400786991225
$ws.Cells.Item(8,5).Value2 = 0.04045134865384774
$ws.Cells.Item(8,6).Value2 = 8.775689685661026
$ws.Cells.Item(8,7).Value2 = 0.001
$ws.Cells.Item(8,8).Value2 = "'False"

# ---- Row 9 (new row added for 21v vs 24.5v) ----
$ws.Cells.Item(9,1).Value2 = "21v"
$ws.Cells.Item(9,2).Value2 = "24.5v"
$ws.Cells.Item(9,3).Value2 = 0.02153969241574086
$ws.Cells.Item(9,4).Value2 = 0.008052603749639412
$ws.Cells.Item(9,5).Value2 = 0.03502678108184232
$ws.Cells.Item(9,6).Value2 = 7.007564632617195
$ws.Cells.Item(9,7).Value2 = 0.001
$ws.Cells.Item(9,8).Value2 = "'False"
